$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.056.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.86%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.180.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.14%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.64%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.69%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.178.07"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.08%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.76%  "

# Row 10
$ws.Range("E10").Value = "  +6.16%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.26"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.20%  "

# Row 12
$ws.Range("E12").Value = "  +2.65%  "

# Row 13
$ws.Range("E13").Value = "  +18.29%  "

# Row 14
$ws.Range("E14").Value = "  +6.11%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.703.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.23%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.141.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.91%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.183.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.24%  "

# Row 18
$ws.Range("E18").Value = "  +5.54%  "

# Row 19
$ws.Range("E19").Value = "  +1.50%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "514.79"
$ws.Range("D20").Style = "Normal"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.88%  "

# Row 22
$ws.Range("E22").Value = "  +6.95%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.86%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.83"
$ws.Range("D24").Style = "Normal"

# Row 25
$ws.Range("E25").Value = "  +3.27%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.08%  "

# Row 27
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.85%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.11%  "

# Row 29
$ws.Range("E29").Value = "  +7.44%  "

# Row 30
$ws.Range("E30").Value = "  +6.34%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +13.30%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.30%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.35"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.15%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.82%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.19%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0900"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.35%  "

# Row 38
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.67%  "

# Row 39
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "476.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.47%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.67%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.075.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.88%  "

# Row 43
$ws.Range("E43").Value = "  +1.93%  "

# Row 44
$ws.Range("E44").Value = "  +5.89%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.47%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.47%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0617"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +19.88%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.19%  "

# Row 51
$ws.Range("E51").Value = "  +2.36%  "
